$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:C3").EntireRow.Delete() | Out-Null

$ws.Range("A21").Value = 0.0018325957935303
$ws.Range("B21").Value = 0.0178678091615438
$ws.Range("C21").Value = 0.0360410511493682
